# Updated cryptos list on Fri Aug  9 15:15:44 UTC 2024 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking strings (e.g. "504.27",
# "1.00", "60.089.21") that must stay stored as plain text, exactly as
# authored - letting Excel auto-convert them to numbers would strip
# trailing zeros / introduce float rounding. Force text storage by
# stamping NumberFormat "@" before assigning, then restore the cell's
# style to match its undecorated neighbour so no stray style index is
# left behind in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell {
    param($row, $value)
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $ws.Range("B$row").Style
}

function Set-PctCell {
    param($row, $value)
    $ws.Range("E$row").Value = $value
}

function Set-CoinCell {
    param($row, $value)
    $ws.Range("B$row").Value = $value
}

function Set-LinkCell {
    param($row, $value)
    $ws.Range("C$row").Value = $value
}

# --- rows 2-39: price / volume(1h) updates ---------------------------------

Set-PriceCell 2  "60.089.21"
Set-PctCell   2  "  +2.77%  "

Set-PriceCell 3  "2.570.14"
Set-PctCell   3  "  +4.47%  "

Set-PriceCell 4  "1.00"
Set-PctCell   4  "  -0.10%  "

Set-PriceCell 5  "504.27"
Set-PctCell   5  "  +1.78%  "

Set-PriceCell 6  "153.40"
Set-PctCell   6  "  -3.62%  "

Set-PriceCell 7  "0.998"
Set-PctCell   7  "  +0.42%  "

Set-PriceCell 8  "0.580"
Set-PctCell   8  "  -5.55%  "

Set-PriceCell 9  "2.593.25"
Set-PctCell   9  "  +3.52%  "

Set-PriceCell 10 "6.60"
Set-PctCell   10 "  +4.51%  "

Set-PctCell   11 "  +0.36%  "

Set-PctCell   12 "  +1.45%  "

Set-PctCell   13 "  +0.62%  "

Set-PriceCell 14 "3.023.04"
Set-PctCell   14 "  +4.61%  "

Set-PriceCell 15 "60.171.07"
Set-PctCell   15 "  +2.82%  "

Set-PriceCell 16 "21.56"
Set-PctCell   16 "  +0.79%  "

Set-PctCell   17 "  +2.54%  "

Set-PriceCell 18 "2.586.50"
Set-PctCell   18 "  +3.95%  "

Set-PriceCell 19 "4.79"
Set-PctCell   19 "  +1.20%  "

Set-PriceCell 20 "346.55"
Set-PctCell   20 "  +5.12%  "

Set-PriceCell 21 "10.28"
Set-PctCell   21 "  +1.15%  "

Set-PriceCell 22 "6.03"
Set-PctCell   22 "  +0.59%  "

Set-PctCell   23 "  +0.02%  "

Set-PriceCell 24 "60.15"
Set-PctCell   24 "  +2.26%  "

Set-PriceCell 25 "0.419"
Set-PctCell   25 "  +2.25%  "

Set-PriceCell 26 "0.166"
Set-PctCell   26 "  +0.58%  "

Set-PriceCell 27 "2.684.58"
Set-PctCell   27 "  +4.31%  "

Set-PriceCell 28 "0.993"
Set-PctCell   28 "  +0.13%  "

Set-PriceCell 29 "0.0₃0845"
Set-PctCell   29 "  +4.35%  "

Set-PctCell   30 "  +0.12%  "

Set-PctCell   31 "  +0.24%  "

Set-PriceCell 32 "155.20"
Set-PctCell   32 "  +2.52%  "

Set-PctCell   33 "  +0.57%  "

Set-PriceCell 34 "1.55"
Set-PctCell   34 "  +0.74%  "

Set-PriceCell 35 "5.70"
Set-PctCell   35 "  +4.60%  "

Set-PctCell   36 "  +3.32%  "

Set-PctCell   37 "  +1.80%  "

Set-PriceCell 38 "0.852"
Set-PctCell   38 "  +20.90%  "

Set-PriceCell 39 "0.841"
Set-PctCell   39 "  -1.23%  "

# --- rows 40/41: Stacks <-> Filecoin swapped, then values refreshed --------

Set-CoinCell  40 "Filecoin"
Set-LinkCell  40 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceCell 40 "3.76"
Set-PctCell   40 "  +3.06%  "

Set-CoinCell  41 "Stacks"
Set-LinkCell  41 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-PriceCell 41 "1.46"
Set-PctCell   41 "  +3.26%  "

# --- rows 42-44 ------------------------------------------------------------

Set-PriceCell 42 "299.55"
Set-PctCell   42 "  +5.67%  "

Set-PctCell   43 "  +3.20%  "

Set-PriceCell 44 "0.0564"
Set-PctCell   44 "  +3.19%  "

# --- rows 45/46: Stellar <-> Mantle swapped, then values refreshed --------

Set-CoinCell  45 "Mantle"
Set-LinkCell  45 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceCell 45 "0.616"
Set-PctCell   45 "  +0.44%  "

Set-CoinCell  46 "Stellar"
Set-LinkCell  46 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceCell 46 "0.0996"
Set-PctCell   46 "  -1.45%  "

# --- rows 47-49 -------------------------------------------------------------

Set-PriceCell 47 "0.998"
Set-PctCell   47 "  +0.41%  "

Set-PriceCell 48 "19.68"
Set-PctCell   48 "  +7.41%  "

Set-PriceCell 49 "4.92"
Set-PctCell   49 "  +3.50%  "

# --- rows 50/51: VeChain <-> Maker swapped, then values refreshed ---------

Set-CoinCell  50 "Maker"
Set-LinkCell  50 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-PriceCell 50 "2.032.21"
Set-PctCell   50 "  +5.83%  "

Set-CoinCell  51 "VeChain"
Set-LinkCell  51 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-PriceCell 51 "0.0233"
Set-PctCell   51 "  -0.80%  "
